$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 583.75
$ws.Range("J2").Value = 749.5
$ws.Range("L2").Value = 749.5
$ws.Range("N2").Value = -975.5
$ws.Range("H33").Value = 11111270
$ws.Range("I33").Value = 12500167
$ws.Range("J33").Value = 90
$ws.Range("K33").Value = 12500167
$ws.Range("L33").Value = 90
$ws.Range("M33").Value = -12499938
$ws.Range("N33").Value = -548
$ws.Range("H64").Value = 12711.375
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 12711.375
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H76").Value = 5749.1665
$ws.Range("I76").Value = 3995
$ws.Range("J76").Value = 6100
$ws.Range("K76").Value = 3995
$ws.Range("L76").Value = 6100
$ws.Range("M76").Value = -3680
$ws.Range("N76").Value = -6730
$ws.Range("H79").Value = 5749.1665
$ws.Range("I79").Value = 3995
$ws.Range("J79").Value = 6100
$ws.Range("K79").Value = 3995
$ws.Range("L79").Value = 6100
$ws.Range("M79").Value = -2903
$ws.Range("N79").Value = -8284
$ws.Range("H92").Value = 10527499
$ws.Range("I92").Value = 15385568
$ws.Range("K92").Value = 15385568
$ws.Range("M92").Value = -15384320
$ws.Range("H100").Value = 1799.75
$ws.Range("J100").Value = 712.5
$ws.Range("L100").Value = 712.5
$ws.Range("N100").Value = -1794.5
$ws.Range("H129").Value = 4065.6667
$ws.Range("I129").Value = 4065.6667
$ws.Range("K129").Value = 12197.0001
$ws.Range("M129").Value = -7197.000100000001
$ws.Range("H132").Value = 9835566
$ws.Range("I132").Value = 11906312
$ws.Range("J132").Value = 172083.33
$ws.Range("K132").Value = 35718936
$ws.Range("L132").Value = 516249.99
$ws.Range("M132").Value = -35716406
$ws.Range("N132").Value = -521309.99
$ws.Range("H133").Value = 196359.1
$ws.Range("J133").Value = 196359.1
$ws.Range("L133").Value = 196359.1
$ws.Range("N133").Value = -206479.1
$ws.Range("H137").Value = 5545.75
$ws.Range("I137").Value = 1276.4445
$ws.Range("K137").Value = 3829.3335
$ws.Range("M137").Value = -1279.3335
$ws.Range("H138").Value = 2932.5671
$ws.Range("I138").Value = 1211.4736
$ws.Range("J138").Value = 3613.8333
$ws.Range("K138").Value = 3634.4208
$ws.Range("L138").Value = 10841.4999
$ws.Range("M138").Value = 1505.5792
$ws.Range("N138").Value = -21121.4999
$ws.Range("H141").Value = 1919
$ws.Range("I141").Value = 1919
$ws.Range("K141").Value = 5757
$ws.Range("M141").Value = -577

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3335.989
$ws.Range("I32").Value = 1713.525
$ws.Range("K32").Value = 1713.525
$ws.Range("M32").Value = -1426.525
$ws.Range("H45").Value = 159707
$ws.Range("I45").Value = 186008.17
$ws.Range("J45").Value = 1900
$ws.Range("K45").Value = 186008.17
$ws.Range("L45").Value = 1900
$ws.Range("M45").Value = -185631.17
$ws.Range("N45").Value = -2654
$ws.Range("H61").Value = 3702.8462
$ws.Range("I61").Value = 2718.0625
$ws.Range("K61").Value = 2718.0625
$ws.Range("M61").Value = -2506.0625
$ws.Range("H74").Value = 34503.555
$ws.Range("I74").Value = 44884.22
$ws.Range("J74").Value = 16137.77
$ws.Range("K74").Value = 44884.22
$ws.Range("L74").Value = 16137.77
$ws.Range("M74").Value = -44010.22
$ws.Range("N74").Value = -17885.77
$ws.Range("H77").Value = 34503.555
$ws.Range("I77").Value = 44884.22
$ws.Range("J77").Value = 16137.77
$ws.Range("K77").Value = 224421.1
$ws.Range("L77").Value = 80688.85000000001
$ws.Range("M77").Value = -220053.1
$ws.Range("N77").Value = -89424.85000000001
$ws.Range("H102").Value = 5301.5713
$ws.Range("I102").Value = 5301.5713
$ws.Range("K102").Value = 5301.5713
$ws.Range("M102").Value = -3679.5713
$ws.Range("H132").Value = 3021.1042
$ws.Range("I132").Value = 2853.4146
$ws.Range("J132").Value = 4003.2856
$ws.Range("K132").Value = 8560.2438
$ws.Range("L132").Value = 12009.8568
$ws.Range("M132").Value = -6030.2438
$ws.Range("N132").Value = -17069.8568
$ws.Range("H136").Value = 3702.8462
$ws.Range("I136").Value = 2718.0625
$ws.Range("K136").Value = 8154.1875
$ws.Range("M136").Value = -5604.1875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4813.857
$ws.Range("I107").Value = 3999.4
$ws.Range("J107").Value = 6850
$ws.Range("K107").Value = 3999.4
$ws.Range("L107").Value = 6850
$ws.Range("M107").Value = -2079.4
$ws.Range("N107").Value = -10690

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25262.756
$ws.Range("I31").Value = 48122.184
$ws.Range("K31").Value = 48122.184
$ws.Range("M31").Value = -47827.184
$ws.Range("H34").Value = 25262.756
$ws.Range("I34").Value = 48122.184
$ws.Range("K34").Value = 48122.184
$ws.Range("M34").Value = -47920.184
$ws.Range("H86").Value = 10224.5
$ws.Range("I86").Value = 9200.5
$ws.Range("K86").Value = 9200.5
$ws.Range("M86").Value = -8077.5
$ws.Range("H89").Value = 10224.5
$ws.Range("I89").Value = 9200.5
$ws.Range("K89").Value = 46002.5
$ws.Range("M89").Value = -40386.5
$ws.Range("H132").Value = 111059.586
$ws.Range("I132").Value = 174203.11
$ws.Range("J132").Value = 3344.1765
$ws.Range("K132").Value = 522609.33
$ws.Range("L132").Value = 10032.5295
$ws.Range("M132").Value = -520079.33
$ws.Range("N132").Value = -15092.5295
$ws.Range("H134").Value = 16551.223
$ws.Range("I134").Value = 14409.095
$ws.Range("J134").Value = 27904.5
$ws.Range("K134").Value = 43227.285
$ws.Range("L134").Value = 83713.5
$ws.Range("M134").Value = -40692.285
$ws.Range("N134").Value = -88783.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 33327428
$ws.Range("I4").Value = 34307050
$ws.Range("J4").Value = 1000000
$ws.Range("K4").Value = 102921150
$ws.Range("L4").Value = 3000000
$ws.Range("M4").Value = -102921038
$ws.Range("N4").Value = -3000224
$ws.Range("H25").Value = 380
$ws.Range("J25").Value = 500
$ws.Range("L25").Value = 1500
$ws.Range("N25").Value = -1838
$ws.Range("H30").Value = 380
$ws.Range("J30").Value = 500
$ws.Range("L30").Value = 1500
$ws.Range("N30").Value = -1704
$ws.Range("H42").Value = 2214
$ws.Range("I42").Value = 250
$ws.Range("J42").Value = 2705
$ws.Range("K42").Value = 750
$ws.Range("L42").Value = 8115
$ws.Range("M42").Value = -216
$ws.Range("N42").Value = -9183
$ws.Range("H140").Value = 3628.25
$ws.Range("I140").Value = 3628.25
$ws.Range("K140").Value = 10884.75
$ws.Range("M140").Value = -5704.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 9561.333000000001
$ws.Range("H93").Value = 29600
$ws.Range("J93").Value = 29600
$ws.Range("L93").Value = 29600
$ws.Range("N93").Value = -33344
$ws.Range("H132").Value = 3009.3125
$ws.Range("I132").Value = 2978.7942
$ws.Range("K132").Value = 8936.382599999999
$ws.Range("M132").Value = -6406.382599999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 1007128
$ws.Range("I122").Value = 1672547
$ws.Range("J122").Value = 8999.5
$ws.Range("K122").Value = 5017641
$ws.Range("L122").Value = 26998.5
$ws.Range("M122").Value = -5015191
$ws.Range("N122").Value = -31898.5
$ws.Range("H136").Value = 2716.6897
$ws.Range("I136").Value = 2410.6843
$ws.Range("J136").Value = 3298.1
$ws.Range("K136").Value = 7232.0529
$ws.Range("L136").Value = 9894.299999999999
$ws.Range("M136").Value = -4682.0529
$ws.Range("N136").Value = -14994.3

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 4749.75
$ws.Range("J6").Value = 4749.75
$ws.Range("L6").Value = 4749.75
$ws.Range("N6").Value = -4979.75
$ws.Range("H136").Value = 1567.9286
$ws.Range("I136").Value = 1354.7812
$ws.Range("J136").Value = 2250
$ws.Range("K136").Value = 4064.3436
$ws.Range("L136").Value = 6750
$ws.Range("M136").Value = -1514.3436
$ws.Range("N136").Value = -11850
